$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-15 Wednesday" "2024-05-16 Thursday"

Replace-Text "32×14=448" "73×60=4380"
Replace-Text "24×90=2160" "73×35=2555"
Replace-Text "19×41=779" "34×71=2414"
Replace-Text "42×47=1974" "13×27=351"
Replace-Text "80×38=3040" "31×22=682"
Replace-Text "66×70=4620" "71×41=2911"
Replace-Text "76×18=1368" "96×42=4032"
Replace-Text "25×75=1875" "20×67=1340"
Replace-Text "66×95=6270" "56×48=2688"
Replace-Text "82×61=5002" "25×99=2475"
Replace-Text "29×87=2523" "97×12=1164"
Replace-Text "57×94=5358" "30×65=1950"
Replace-Text "75×75=5625" "60×11=660"
Replace-Text "20×28=560" "33×98=3234"
Replace-Text "65×81=5265" "18×68=1224"
Replace-Text "18×13=234" "87×45=3915"
Replace-Text "27×62=1674" "16×27=432"
Replace-Text "97×49=4753" "83×24=1992"
Replace-Text "86×67=5762" "39×34=1326"
Replace-Text "96×29=2784" "97×90=8730"
Replace-Text "33×11=363" "65×78=5070"
Replace-Text "55×66=3630" "30×18=540"
Replace-Text "87×28=2436" "33×80=2640"
Replace-Text "97×19=1843" "35×20=700"
Replace-Text "49×85=4165" "94×93=8742"

Write-Output "Done"
